# Rework the scenario sheet from "movie catalog" steps to "marketplace" steps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old content (rows 1-5, cols A-E) so only the new 2-row scenario remains.
$ws.Cells.Clear()

# Row 1 - first scenario ("search & buy"), written left to right so the
# shared-string table is populated in the same order as the target file.
$ws.Range("A1").Value = "Захожу в католог"
$ws.Range("B1").Value = "Ввожу в поисковой строке  что хочешь купить"
$ws.Range("C1").Value = "Нажимаешь кнопку поиск"
$ws.Range("F1").Value = "переходишь по ссылке и покупаешь вещь"

# Row 2 - second scenario ("register & sell an item").
$ws.Range("A2").Value = "Захожу в католог"
$ws.Range("B2").Value = "Регестрируешься"
$ws.Range("C2").Value = "Добавляешь свою необычную вещь/продукт"

# D1 is filled last so it reuses the trailing slot in the shared-string table.
$ws.Range("D1").Value = "Сортируешь по: отзывам, стоимости, региону, рейтингу, времени прибытия, типу товара"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("C5").Select()
